# Remove stray ".0" suffixes left over from a float->string conversion on the
# "compte" (account code) column, e.g. "601.0" -> "601".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and ($val -is [string]) -and ($val -match '^\d+\.0$')) {
            $newVal = $val -replace '\.0$', ''
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
            $cell.Style = "Normal"
        }
    }
}
